$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(2, 6).Value = 25
$ws.Cells.Item(2, 8).Value = 25
$ws.Cells.Item(3, 6).Value = 13
$ws.Cells.Item(3, 8).Value = 13
$ws.Cells.Item(8, 5).Value = 14
$ws.Cells.Item(15, 6).Value = 81
$ws.Cells.Item(15, 8).Value = 81
$ws.Cells.Item(16, 6).Value = 4
$ws.Cells.Item(16, 8).Value = 4
$ws.Cells.Item(17, 6).Value = 48
$ws.Cells.Item(17, 8).Value = 48
$ws.Cells.Item(18, 6).Value = 45
$ws.Cells.Item(18, 8).Value = 45
$ws.Cells.Item(19, 6).Value = 25
$ws.Cells.Item(19, 8).Value = 25
$ws.Cells.Item(20, 5).Value = 5
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(20, 8).Value = 2
$ws.Cells.Item(25, 5).Value = 20
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(27, 8).Value = 4
$ws.Cells.Item(28, 6).Value = 10
$ws.Cells.Item(28, 8).Value = 10
$ws.Cells.Item(39, 6).Value = 14
$ws.Cells.Item(39, 8).Value = 14
$ws.Cells.Item(41, 5).Value = 37
$ws.Cells.Item(42, 6).Value = 15
$ws.Cells.Item(42, 8).Value = 15
$ws.Cells.Item(44, 5).Value = 27
$ws.Cells.Item(44, 6).Value = 13
$ws.Cells.Item(44, 8).Value = 13
$ws.Cells.Item(47, 6).Value = 34
$ws.Cells.Item(47, 8).Value = 34
$ws.Cells.Item(55, 6).Value = 3
$ws.Cells.Item(55, 8).Value = 3
$ws.Cells.Item(57, 5).Value = 13
$ws.Cells.Item(66, 6).Value = 20
$ws.Cells.Item(66, 8).Value = 20
$ws.Cells.Item(72, 5).Value = 37
$ws.Cells.Item(74, 6).Value = 6
$ws.Cells.Item(74, 8).Value = 6
$ws.Cells.Item(78, 5).Value = 42
$ws.Cells.Item(79, 5).Value = 32
$ws.Cells.Item(89, 6).Value = 14
$ws.Cells.Item(89, 8).Value = 14
